$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.867.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.28%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.41%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4433"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.59%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3669"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.82%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07349"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.41%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8566"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.62%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.75%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.806.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.599"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.79%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.46"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.95%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.304"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07066"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.47%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008727"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.42%  "

# Row 19
$ws.Range("E19").Value = "  -0.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.40%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.887.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.47%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.152"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.77%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.52%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.25%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.176"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.22%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.212"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.40%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.27%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08836"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.7509"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.174"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.59%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.929"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.16%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.453"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.07%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.43%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.084"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.96%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01965"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05201"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.03%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5312"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.30%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.863"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.45%  "

# Row 41
$ws.Range("E41").Value = "  -4.32%  "

# Row 42
$ws.Range("E42").Value = "  -0.34%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5161"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.78%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.418"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.56%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.988"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.33%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "105.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.87%  "

# Row 48
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.45%  "

# Row 49
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.667"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06320"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.08%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9193"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.54%  "
